# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Bad Drivers section (row 3 = Intel AX201 22.120.0.3 entry, row 4 = Totals row)
$ws.Range("C3").Value = 337
$ws.Range("D3").Value = 95.40000000000001
$ws.Range("C4").Value = 337

# Good Drivers section (row 12 = AX201 23.100.0.4, row 13 = AX201 22.80.0.9)
$ws.Range("B12").Value = 449371
$ws.Range("B13").Value = 77999
